$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.7
$ws.Range("D2").Value = 162
$ws.Range("E2").Value = 33
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 71
$ws.Range("H2").Value = 0.92
$ws.Range("I2").Value = 0.7
$ws.Range("J2").Value = 0.8
$ws.Range("K2").Value = 0.7
$ws.Range("L2").Value = 0.7

$ws.Range("B3").Value = 0.87
$ws.Range("D3").Value = 192
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 22
$ws.Range("G3").Value = 41
$ws.Range("H3").Value = 0.9
$ws.Range("I3").Value = 0.82
$ws.Range("J3").Value = 0.86
$ws.Range("K3").Value = 0.53
$ws.Range("L3").Value = 0.68

$ws.Range("B4").Value = 0.84
$ws.Range("E4").Value = 32
$ws.Range("F4").Value = 15
$ws.Range("K4").Value = 0.68
$ws.Range("L4").Value = 0.6899999999999999

$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 0.85
$ws.Range("D5").Value = 216
$ws.Range("E5").Value = 21
$ws.Range("F5").Value = 26
$ws.Range("G5").Value = 17
$ws.Range("H5").Value = 0.89
$ws.Range("I5").Value = 0.93
$ws.Range("J5").Value = 0.91
$ws.Range("K5").Value = 0.45

$ws.Range("B6").Value = 0.97
$ws.Range("C6").Value = 0.76
$ws.Range("D6").Value = 233
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 47
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0.83
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 0.91
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0.5
